# BOM.xlsx update:
#  - Added LED matrix spacing (D45 LED Grün, D46 LED Rot)
#  - Added Soft-Reset / Hard-Reset buttons (SW1-SW3)
#  - Added USB-Interface (FT232RL + USB Type B connector)
#  - Added Programming-Button
#  - Added ISP-Programmer
#  - Created 3D Frontplate mockup
#
# Appends rows 12-18 to the Bill-of-Materials table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - FT232RL USB/Serial bridge
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "U"
$ws.Range("C12").Value = "FT232RL"
$ws.Range("D12").Value = "U4"
$ws.Range("E12").Value = "FT 232 RL"

# Row 13 - USB Type B connector
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "J"
$ws.Range("C13").Value = "USB Type B"
$ws.Range("D13").Value = "J1"
$ws.Range("E13").Value = "RND 205-00858"

# Row 14 - Green LED
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "D"
$ws.Range("C14").Value = "LED Grün"
$ws.Range("D14").Value = "D45"
$ws.Range("E14").Value = "EVL 17-21SYGC/S2"

# Row 15 - Red LED
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "D"
$ws.Range("C15").Value = "LED Rot"
$ws.Range("D15").Value = "D46"
$ws.Range("E15").Value = "KBT KP-2012EC"

# Row 16 - 220 Ohm resistors
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "R"
$ws.Range("C16").Value = 220
$ws.Range("D16").Value = "R43, R44"
$ws.Range("E16").Value = "RND 0805 1 220"

# Row 17 - 470 Ohm resistor
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "R"
$ws.Range("C17").Value = 470
$ws.Range("D17").Value = "R45"

# Row 18 - SMD push buttons (Soft-Reset / Hard-Reset / Programming)
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "SW"
$ws.Range("C18").Value = "SMD Taster"
$ws.Range("D18").Value = "SW1 - SW3"
$ws.Range("E18").Value = "TASTER 1612.11"

# Match the author's final selection state
$ws.Range("E18").Select()
